$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.974.26"
$ws.Range("E2").Value = "  -1.33%  "
$ws.Range("D3").Value = "2.302.73"
$ws.Range("E3").Value = "  -1.64%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "300.80"
$ws.Range("E5").Value = "  -1.37%  "
$ws.Range("D6").Value = "96.94"
$ws.Range("E6").Value = "  -4.93%  "
$ws.Range("D7").Value = "0.506"
$ws.Range("E7").Value = "  -1.30%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "0.496"
$ws.Range("E9").Value = "  -3.54%  "
$ws.Range("D10").Value = "33.65"
$ws.Range("E10").Value = "  -4.67%  "
$ws.Range("D11").Value = "0.0793"
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("D12").Value = "49.26"
$ws.Range("E12").Value = "  -4.63%  "
$ws.Range("E13").Value = "  +2.15%  "
$ws.Range("D14").Value = "16.83"
$ws.Range("E14").Value = "  +7.77%  "
$ws.Range("E15").Value = "  -0.52%  "
$ws.Range("D16").Value = "2.663.19"
$ws.Range("E16").Value = "  -1.49%  "
$ws.Range("D17").Value = "2.324.34"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").Value = "0.807"
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("D19").Value = "42.907.80"
$ws.Range("E19").Value = "  -1.29%  "
$ws.Range("D20").Value = "11.60"
$ws.Range("E20").Value = "  -2.19%  "
$ws.Range("D21").Value = "0.0₃0900"
$ws.Range("E21").Value = "  -1.06%  "
$ws.Range("D22").Value = "6.00"
$ws.Range("E22").Value = "  -2.19%  "
$ws.Range("D23").Value = "67.15"
$ws.Range("E23").Value = "  -1.79%  "
$ws.Range("D24").Value = "236.32"
$ws.Range("E24").Value = "  -0.87%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").Value = "2.45"
$ws.Range("E27").Value = "  -3.54%  "
$ws.Range("D28").Value = "24.87"
$ws.Range("E28").Value = "  -0.80%  "
$ws.Range("E29").Value = "  +3.61%  "
$ws.Range("D30").Value = "166.71"
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("D31").Value = "33.79"
$ws.Range("E31").Value = "  -2.56%  "
$ws.Range("E32").Value = "  -1.51%  "
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("D34").Value = "4.78"
$ws.Range("E34").Value = "  +5.87%  "
$ws.Range("E35").Value = "  -2.41%  "
$ws.Range("E36").Value = "  -0.82%  "
$ws.Range("D37").Value = "16.83"
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("D38").Value = "0.0695"
$ws.Range("E38").Value = "  -1.72%  "
$ws.Range("E39").Value = "  -3.48%  "
$ws.Range("E40").Value = "  -2.15%  "
$ws.Range("E41").Value = "  -4.49%  "
$ws.Range("E42").Value = "  -2.14%  "
$ws.Range("E43").Value = "  -2.61%  "
$ws.Range("D44").Value = "1.971.83"
$ws.Range("E44").Value = "  -1.16%  "
$ws.Range("D45").Value = "0.0280"
$ws.Range("E45").Value = "  -2.01%  "
$ws.Range("D46").Value = "17.60"
$ws.Range("E46").Value = "  -4.91%  "
$ws.Range("D47").Value = "9.79"
$ws.Range("E47").Value = "  -2.25%  "
$ws.Range("E48").Value = "  -3.66%  "
$ws.Range("D49").Value = "2.528.75"
$ws.Range("E49").Value = "  -1.44%  "
$ws.Range("D50").Value = "52.69"
$ws.Range("E50").Value = "  -7.47%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "1.50"
$ws.Range("E51").Value = "  -3.15%  "
